# Excel COM-interop edit script
# Updates the currentAveragePrice / LevePrice / LeveProfit columns (H,I,J,K,L,M,N)
# for specific Leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets, matching a
# scheduled market-price refresh run. A couple of rows also gain/lose a previously
# empty M/N cell where the recomputed profit crossed from blank to a real number.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2: Mercury Rising
$ws.Range("H2").Value = 1589.6666
$ws.Range("I2").Value = 1589.6666
$ws.Range("K2").Value = 1589.6666
$ws.Range("M2").Value = -1476.6666
# Row 28: The Writing Is Not on the Wall
$ws.Range("H28").Value = 925.5454999999999
$ws.Range("J28").Value = 1331.75
$ws.Range("L28").Value = 1331.75
$ws.Range("N28").Value = -2301.75
# Row 99: Rumor Has It
$ws.Range("H99").Value = 2799.4614
$ws.Range("J99").Value = 4054.75
$ws.Range("L99").Value = 12164.25
$ws.Range("N99").Value = -15160.25
# Row 101: Edge of the Arcane
$ws.Range("H101").Value = 1415
$ws.Range("I101").Value = 1223.5714
$ws.Range("J101").Value = 4095
$ws.Range("K101").Value = 3670.7142
$ws.Range("L101").Value = 12285
$ws.Range("M101").Value = -2048.7142
$ws.Range("N101").Value = -15529
# Row 107: Another Man's Ink
$ws.Range("H107").Value = 1597.8182
$ws.Range("I107").Value = 1662.25
$ws.Range("J107").Value = 1426
$ws.Range("K107").Value = 1662.25
$ws.Range("L107").Value = 1426
$ws.Range("M107").Value = 257.75
$ws.Range("N107").Value = -5266
# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1203.6875
$ws.Range("I137").Value = 1271.5834
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 3814.7502
$ws.Range("L137").Value = 3000
$ws.Range("M137").Value = -1264.7502
$ws.Range("N137").Value = -8100
# Row 138: All-night Crafting
$ws.Range("H138").Value = 6805.5884
$ws.Range("I138").Value = 3735
$ws.Range("J138").Value = 8085
$ws.Range("K138").Value = 11205
$ws.Range("L138").Value = 24255
$ws.Range("M138").Value = -6065
$ws.Range("N138").Value = -34535
# Row 141: Remedy for Reason
$ws.Range("H141").Value = 3213
$ws.Range("I141").Value = 3042.875
$ws.Range("K141").Value = 9128.625
$ws.Range("M141").Value = -3948.625

$ws = $wb.Worksheets.Item("ARM")
# Row 4: Eyes Bigger than the Plate
$ws.Range("H4").Value = 198.5
$ws.Range("J4").Value = 199.66667
$ws.Range("L4").Value = 199.66667
$ws.Range("N4").Value = -431.66667
# Row 29: No Hand-me-downs
$ws.Range("H29").Value = 15666.333
$ws.Range("J29").Value = 15666.333
$ws.Range("L29").Value = 15666.333
$ws.Range("N29").Value = -16282.333
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 3742.0833
$ws.Range("I32").Value = 3823.3044
$ws.Range("K32").Value = 3823.3044
$ws.Range("M32").Value = -3536.3044
# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 1928.762
$ws.Range("I45").Value = 1714.4667
$ws.Range("J45").Value = 2464.5
$ws.Range("K45").Value = 1714.4667
$ws.Range("L45").Value = 2464.5
$ws.Range("M45").Value = -1337.4667
$ws.Range("N45").Value = -3218.5
# Row 50: Liquid Persuasion
$ws.Range("H50").Value = 10981.637
$ws.Range("I50").Value = 4933.1665
$ws.Range("J50").Value = 18239.8
$ws.Range("K50").Value = 4933.1665
$ws.Range("L50").Value = 18239.8
$ws.Range("M50").Value = -4219.1665
$ws.Range("N50").Value = -19667.8
# Row 63: Rivets Run through It
$ws.Range("H63").Value = 1895.8572
$ws.Range("J63").Value = 1823.625
$ws.Range("L63").Value = 1823.625
$ws.Range("N63").Value = -3195.625
# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 1895.8572
$ws.Range("J66").Value = 1823.625
$ws.Range("L66").Value = 9118.125
$ws.Range("N66").Value = -15982.125
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 774929.5600000001
$ws.Range("I74").Value = 3619.6453
$ws.Range("K74").Value = 3619.6453
$ws.Range("M74").Value = -2745.6453
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 774929.5600000001
$ws.Range("I77").Value = 3619.6453
$ws.Range("K77").Value = 18098.2265
$ws.Range("M77").Value = -13730.2265
# Row 109: A Head of Demand
$ws.Range("H109").Value = 44999
$ws.Range("J109").Value = 44999
$ws.Range("L109").Value = 44999
$ws.Range("N109").Value = -47773
# Row 122: Haste for High Durium
$ws.Range("H122").Value = 3185.75
$ws.Range("I122").Value = 2965.0476
$ws.Range("J122").Value = 4730.6665
$ws.Range("K122").Value = 8895.1428
$ws.Range("L122").Value = 14191.9995
$ws.Range("M122").Value = -6445.1428
$ws.Range("N122").Value = -19091.9995
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 45463924
$ws.Range("J132").Value = 125018740
$ws.Range("L132").Value = 375056220
$ws.Range("N132").Value = -375061280

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt
$ws.Range("H20").Value = 4073.3076
$ws.Range("I20").Value = 3550.7144
$ws.Range("K20").Value = 3550.7144
$ws.Range("M20").Value = -3303.7144
# Row 99: Meddle in Metal
$ws.Range("H99").Value = 1900
$ws.Range("I99").Value = 1900
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1900
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -402
$ws.Range("N99").ClearContents()
# Row 110: Selective Logging
$ws.Range("H110").Value = 44686.75
$ws.Range("J110").Value = 44686.75
$ws.Range("L110").Value = 44686.75
$ws.Range("N110").Value = -52866.75

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 5528
$ws.Range("I31").Value = 8500
$ws.Range("K31").Value = 8500
$ws.Range("M31").Value = -8205
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 5528
$ws.Range("I34").Value = 8500
$ws.Range("K34").Value = 8500
$ws.Range("M34").Value = -8298

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Range("H5").Value = 1939.3
$ws.Range("I5").Value = 1233
$ws.Range("K5").Value = 3699
$ws.Range("M5").Value = -3587
# Row 109: Cure for What Ails
$ws.Range("H109").Value = 5471.1665
$ws.Range("I109").Value = 4565.4
$ws.Range("K109").Value = 13696.2
$ws.Range("M109").Value = -12656.2
# Row 115: Mixology
$ws.Range("H115").Value = 6400
$ws.Range("I115").Value = 1000
$ws.Range("K115").Value = 3000
$ws.Range("M115").Value = -1825
# Row 120: A Happy End
$ws.Range("H120").Value = 30000
$ws.Range("I120").Value = 30000
$ws.Range("K120").Value = 90000
$ws.Range("M120").Value = -85162
# Row 132: More Mezcal
$ws.Range("H132").Value = 1766
$ws.Range("I132").Value = 1124.5
$ws.Range("J132").Value = 2279.2
$ws.Range("K132").Value = 10120.5
$ws.Range("L132").Value = 20512.8
$ws.Range("M132").Value = -7590.5
$ws.Range("N132").Value = -25572.8
# Row 133: Friends Are Food
$ws.Range("H133").Value = 9676.666999999999
$ws.Range("I133").Value = 9030
$ws.Range("K133").Value = 27090
$ws.Range("M133").Value = -22030
# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 1939.3
$ws.Range("I135").Value = 1233
$ws.Range("K135").Value = 11097
$ws.Range("M135").Value = -8562

$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 1557.0435
$ws.Range("I97").Value = 1434.5
$ws.Range("K97").Value = 1434.5
$ws.Range("M97").Value = -938.5
# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 2935.375
$ws.Range("I113").Value = 3047.8
$ws.Range("J113").Value = 2748
$ws.Range("K113").Value = 3047.8
$ws.Range("L113").Value = 2748
$ws.Range("M113").Value = -877.8000000000002
$ws.Range("N113").Value = -7088

$ws = $wb.Worksheets.Item("LTW")
# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 1098.6666
$ws.Range("I93").Value = 1118.4
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 1118.4
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = 129.5999999999999
$ws.Range("N93").Value = -3496
